$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.8781003333333334
$ws.Range("H2").Value = 2.634301
$ws.Range("I2").Value = 0.1010434633250494
$ws.Range("J2").Value = 0.1010434633250494
$ws.Range("M2").Value = 42.04602466666667
$ws.Range("N2").Value = 126.138074
$ws.Range("O2").Value = 0.1180439555498783
$ws.Range("P2").Value = 0.1180439555498783
$ws.Range("Q2").Value = 36.92062827514156
$ws.Range("R2").Value = 332.285654476274
$ws.Range("S2").Value = 0.01192757009334789
$ws.Range("T2").Value = 0.01192757009334789

$ws.Range("G3").Value = 0.8781003333333334
$ws.Range("H3").Value = 2.634301
$ws.Range("I3").Value = 0.1010434633250494
$ws.Range("J3").Value = 0.1010434633250494
$ws.Range("O3").Value = 0.1482760805823429
$ws.Range("P3").Value = 0.1482760805823429
$ws.Range("Q3").Value = 46.37633521999734
$ws.Range("R3").Value = 417.387016979976
$ws.Range("S3").Value = 0.01498232871030404
$ws.Range("T3").Value = 0.01498232871030404

$ws.Range("G4").Value = 0.8781003333333334
$ws.Range("H4").Value = 2.634301
$ws.Range("I4").Value = 0.1010434633250494
$ws.Range("J4").Value = 0.1010434633250494
$ws.Range("M4").Value = 68.81807333333334
$ws.Range("N4").Value = 206.45422
$ws.Range("O4").Value = 0.1932063174578422
$ws.Range("P4").Value = 0.1932063174578422
$ws.Range("Q4").Value = 60.42917313335779
$ws.Range("R4").Value = 543.8625582002201
$ws.Range("S4").Value = 0.01952223545221933
$ws.Range("T4").Value = 0.01952223545221933

$ws.Range("G5").Value = 0.8781003333333334
$ws.Range("H5").Value = 2.634301
$ws.Range("I5").Value = 0.1010434633250494
$ws.Range("J5").Value = 0.1010434633250494
$ws.Range("M5").Value = 11.78107633333333
$ws.Range("N5").Value = 35.343229
$ws.Range("O5").Value = 0.03307529931894448
$ws.Range("P5").Value = 0.03307529931894448
$ws.Range("Q5").Value = 10.34496705532544
$ws.Range("R5").Value = 93.10470349792901
$ws.Range("S5").Value = 0.003342042793698798
$ws.Range("T5").Value = 0.003342042793698798

$ws.Range("G6").Value = 0.8781003333333334
$ws.Range("H6").Value = 2.634301
$ws.Range("I6").Value = 0.1010434633250494
$ws.Range("J6").Value = 0.1010434633250494
$ws.Range("M6").Value = 139.820737
$ws.Range("N6").Value = 419.462211
$ws.Range("O6").Value = 0.3925458588351179
$ws.Range("P6").Value = 0.3925458588351179
$ws.Range("Q6").Value = 122.7766357666123
$ws.Range("R6").Value = 1104.989721899511
$ws.Range("S6").Value = 0.03966419309060626
$ws.Range("T6").Value = 0.03966419309060626

$ws.Range("G7").Value = 0.8781003333333334
$ws.Range("H7").Value = 2.634301
$ws.Range("I7").Value = 0.1010434633250494
$ws.Range("J7").Value = 0.1010434633250494
$ws.Range("M7").Value = 40.90925733333334
$ws.Range("N7").Value = 122.727772
$ws.Range("O7").Value = 0.1148524882558742
$ws.Range("P7").Value = 0.1148524882558742
$ws.Range("Q7").Value = 35.92243250081911
$ws.Range("R7").Value = 323.301892507372
$ws.Range("S7").Value = 0.0116050931848731
$ws.Range("T7").Value = 0.01160509318487309

$ws.Range("I8").Value = 0.01326751606355713
$ws.Range("J8").Value = 0.01326751606355713
$ws.Range("M8").Value = 42.04602466666667
$ws.Range("N8").Value = 126.138074
$ws.Range("O8").Value = 0.1180439555498783
$ws.Range("P8").Value = 0.1180439555498783
$ws.Range("Q8").Value = 4.847864598042
$ws.Range("R8").Value = 43.630781382378
$ws.Range("S8").Value = 0.001566150076463834
$ws.Range("T8").Value = 0.001566150076463834

$ws.Range("I9").Value = 0.01326751606355713
$ws.Range("J9").Value = 0.01326751606355713
$ws.Range("O9").Value = 0.1482760805823429
$ws.Range("P9").Value = 0.1482760805823429
$ws.Range("S9").Value = 0.001967255280967526
$ws.Range("T9").Value = 0.001967255280967526

$ws.Range("I10").Value = 0.01326751606355713
$ws.Range("J10").Value = 0.01326751606355713
$ws.Range("M10").Value = 68.81807333333334
$ws.Range("N10").Value = 206.45422
$ws.Range("O10").Value = 0.1932063174578422
$ws.Range("P10").Value = 0.1932063174578422
$ws.Range("Q10").Value = 7.934655037260002
$ws.Range("R10").Value = 71.41189533534001
$ws.Range("S10").Value = 0.00256336792045264
$ws.Range("T10").Value = 0.00256336792045264

$ws.Range("I11").Value = 0.01326751606355713
$ws.Range("J11").Value = 0.01326751606355713
$ws.Range("M11").Value = 11.78107633333333
$ws.Range("N11").Value = 35.343229
$ws.Range("O11").Value = 0.03307529931894448
$ws.Range("P11").Value = 0.03307529931894448
$ws.Range("Q11").Value = 1.358346320157
$ws.Range("R11").Value = 12.225116881413
$ws.Range("S11").Value = 0.0004388270650210561
$ws.Range("T11").Value = 0.0004388270650210561

$ws.Range("I12").Value = 0.01326751606355713
$ws.Range("J12").Value = 0.01326751606355713
$ws.Range("M12").Value = 139.820737
$ws.Range("N12").Value = 419.462211
$ws.Range("O12").Value = 0.3925458588351179
$ws.Range("P12").Value = 0.3925458588351179
$ws.Range("Q12").Value = 16.121191155363
$ws.Range("R12").Value = 145.090720398267
$ws.Range("S12").Value = 0.005208108487777757
$ws.Range("T12").Value = 0.005208108487777757

$ws.Range("I13").Value = 0.01326751606355713
$ws.Range("J13").Value = 0.01326751606355713
$ws.Range("M13").Value = 40.90925733333334
$ws.Range("N13").Value = 122.727772
$ws.Range("O13").Value = 0.1148524882558742
$ws.Range("P13").Value = 0.1148524882558742
$ws.Range("Q13").Value = 4.716796461276
$ws.Range("R13").Value = 42.451168151484
$ws.Range("S13").Value = 0.001523807232874318
$ws.Range("T13").Value = 0.001523807232874318

$ws.Range("G14").Value = 3.520787
$ws.Range("H14").Value = 10.562361
$ws.Range("I14").Value = 0.4051387963370292
$ws.Range("J14").Value = 0.4051387963370292
$ws.Range("M14").Value = 42.04602466666667
$ws.Range("N14").Value = 126.138074
$ws.Range("O14").Value = 0.1180439555498783
$ws.Range("P14").Value = 0.1180439555498783
$ws.Range("Q14").Value = 148.0350970480793
$ws.Range("R14").Value = 1332.315873432714
$ws.Range("S14").Value = 0.04782418606633945
$ws.Range("T14").Value = 0.04782418606633945

$ws.Range("G15").Value = 3.520787
$ws.Range("H15").Value = 10.562361
$ws.Range("I15").Value = 0.4051387963370292
$ws.Range("J15").Value = 0.4051387963370292
$ws.Range("O15").Value = 0.1482760805823429
$ws.Range("P15").Value = 0.1482760805823429
$ws.Range("Q15").Value = 185.948224766504
$ws.Range("R15").Value = 1673.534022898536
$ws.Range("S15").Value = 0.06007239281270275
$ws.Range("T15").Value = 0.06007239281270275

$ws.Range("G16").Value = 3.520787
$ws.Range("H16").Value = 10.562361
$ws.Range("I16").Value = 0.4051387963370292
$ws.Range("J16").Value = 0.4051387963370292
$ws.Range("M16").Value = 68.81807333333334
$ws.Range("N16").Value = 206.45422
$ws.Range("O16").Value = 0.1932063174578422
$ws.Range("P16").Value = 0.1932063174578422
$ws.Range("Q16").Value = 242.2937779570467
$ws.Range("R16").Value = 2180.64400161342
$ws.Range("S16").Value = 0.07827537489958013
$ws.Range("T16").Value = 0.07827537489958011

$ws.Range("G17").Value = 3.520787
$ws.Range("H17").Value = 10.562361
$ws.Range("I17").Value = 0.4051387963370292
$ws.Range("J17").Value = 0.4051387963370292
$ws.Range("M17").Value = 11.78107633333333
$ws.Range("N17").Value = 35.343229
$ws.Range("O17").Value = 0.03307529931894448
$ws.Range("P17").Value = 0.03307529931894448
$ws.Range("Q17").Value = 41.47866040040766
$ws.Range("R17").Value = 373.307943603669
$ws.Range("S17").Value = 0.01340008695456413
$ws.Range("T17").Value = 0.01340008695456413

$ws.Range("G18").Value = 3.520787
$ws.Range("H18").Value = 10.562361
$ws.Range("I18").Value = 0.4051387963370292
$ws.Range("J18").Value = 0.4051387963370292
$ws.Range("M18").Value = 139.820737
$ws.Range("N18").Value = 419.462211
$ws.Range("O18").Value = 0.3925458588351179
$ws.Range("P18").Value = 0.3925458588351179
$ws.Range("Q18").Value = 492.279033160019
$ws.Range("R18").Value = 4430.511298440171
$ws.Range("S18").Value = 0.159035556755545
$ws.Range("T18").Value = 0.159035556755545

$ws.Range("G19").Value = 3.520787
$ws.Range("H19").Value = 10.562361
$ws.Range("I19").Value = 0.4051387963370292
$ws.Range("J19").Value = 0.4051387963370292
$ws.Range("M19").Value = 40.90925733333334
$ws.Range("N19").Value = 122.727772
$ws.Range("O19").Value = 0.1148524882558742
$ws.Range("P19").Value = 0.1148524882558742
$ws.Range("Q19").Value = 144.0327813988547
$ws.Range("R19").Value = 1296.295032589692
$ws.Range("S19").Value = 0.04653119884829766
$ws.Range("T19").Value = 0.04653119884829766

$ws.Range("G20").Value = 0.0464
$ws.Range("H20").Value = 0.1392
$ws.Range("I20").Value = 0.005339272199663925
$ws.Range("J20").Value = 0.005339272199663925
$ws.Range("M20").Value = 42.04602466666667
$ws.Range("N20").Value = 126.138074
$ws.Range("O20").Value = 0.1180439555498783
$ws.Range("P20").Value = 0.1180439555498783
$ws.Range("Q20").Value = 1.950935544533333
$ws.Range("R20").Value = 17.5584199008
$ws.Range("S20").Value = 0.0006302688102058291
$ws.Range("T20").Value = 0.0006302688102058291

$ws.Range("G21").Value = 0.0464
$ws.Range("H21").Value = 0.1392
$ws.Range("I21").Value = 0.005339272199663925
$ws.Range("J21").Value = 0.005339272199663925
$ws.Range("O21").Value = 0.1482760805823429
$ws.Range("P21").Value = 0.1482760805823429
$ws.Range("Q21").Value = 2.4505877888
$ws.Range("R21").Value = 22.0552900992
$ws.Range("S21").Value = 0.0007916863549284316
$ws.Range("T21").Value = 0.0007916863549284316

$ws.Range("G22").Value = 0.0464
$ws.Range("H22").Value = 0.1392
$ws.Range("I22").Value = 0.005339272199663925
$ws.Range("J22").Value = 0.005339272199663925
$ws.Range("M22").Value = 68.81807333333334
$ws.Range("N22").Value = 206.45422
$ws.Range("O22").Value = 0.1932063174578422
$ws.Range("P22").Value = 0.1932063174578422
$ws.Range("Q22").Value = 3.193158602666667
$ws.Range("R22").Value = 28.738427424
$ws.Range("S22").Value = 0.0010315811196021
$ws.Range("T22").Value = 0.0010315811196021

$ws.Range("G23").Value = 0.0464
$ws.Range("H23").Value = 0.1392
$ws.Range("I23").Value = 0.005339272199663925
$ws.Range("J23").Value = 0.005339272199663925
$ws.Range("M23").Value = 11.78107633333333
$ws.Range("N23").Value = 35.343229
$ws.Range("O23").Value = 0.03307529931894448
$ws.Range("P23").Value = 0.03307529931894448
$ws.Range("Q23").Value = 0.5466419418666666
$ws.Range("R23").Value = 4.9197774768
$ws.Range("S23").Value = 0.0001765980261492034
$ws.Range("T23").Value = 0.0001765980261492034

$ws.Range("G24").Value = 0.0464
$ws.Range("H24").Value = 0.1392
$ws.Range("I24").Value = 0.005339272199663925
$ws.Range("J24").Value = 0.005339272199663925
$ws.Range("M24").Value = 139.820737
$ws.Range("N24").Value = 419.462211
$ws.Range("O24").Value = 0.3925458588351179
$ws.Range("P24").Value = 0.3925458588351179
$ws.Range("Q24").Value = 6.4876821968
$ws.Range("R24").Value = 58.3891397712
$ws.Range("S24").Value = 0.002095909191171545
$ws.Range("T24").Value = 0.002095909191171545

$ws.Range("G25").Value = 0.0464
$ws.Range("H25").Value = 0.1392
$ws.Range("I25").Value = 0.005339272199663925
$ws.Range("J25").Value = 0.005339272199663925
$ws.Range("M25").Value = 40.90925733333334
$ws.Range("N25").Value = 122.727772
$ws.Range("O25").Value = 0.1148524882558742
$ws.Range("P25").Value = 0.1148524882558742
$ws.Range("Q25").Value = 1.898189540266667
$ws.Range("R25").Value = 17.0837058624
$ws.Range("S25").Value = 0.0006132286976068167
$ws.Range("T25").Value = 0.0006132286976068167

$ws.Range("E26").Value = 2
$ws.Range("F26").Value = 0.6666666666666666
$ws.Range("G26").Value = 0.41099
$ws.Range("H26").Value = 1.23297
$ws.Range("I26").Value = 0.0472928336495663
$ws.Range("J26").Value = 0.0472928336495663
$ws.Range("M26").Value = 42.04602466666667
$ws.Range("N26").Value = 126.138074
$ws.Range("O26").Value = 0.1180439555498783
$ws.Range("P26").Value = 0.1180439555498783
$ws.Range("Q26").Value = 17.28049567775333
$ws.Range("R26").Value = 155.52446109978
$ws.Range("S26").Value = 0.005582633153157192
$ws.Range("T26").Value = 0.005582633153157192

$ws.Range("E27").Value = 2
$ws.Range("F27").Value = 0.6666666666666666
$ws.Range("G27").Value = 0.41099
$ws.Range("H27").Value = 1.23297
$ws.Range("I27").Value = 0.0472928336495663
$ws.Range("J27").Value = 0.0472928336495663
$ws.Range("O27").Value = 0.1482760805823429
$ws.Range("P27").Value = 0.1482760805823429
$ws.Range("Q27").Value = 21.70618696808
$ws.Range("R27").Value = 195.35568271272
$ws.Range("S27").Value = 0.007012396013190432
$ws.Range("T27").Value = 0.007012396013190432

$ws.Range("E28").Value = 2
$ws.Range("F28").Value = 0.6666666666666666
$ws.Range("G28").Value = 0.41099
$ws.Range("H28").Value = 1.23297
$ws.Range("I28").Value = 0.0472928336495663
$ws.Range("J28").Value = 0.0472928336495663
$ws.Range("M28").Value = 68.81807333333334
$ws.Range("N28").Value = 206.45422
$ws.Range("O28").Value = 0.1932063174578422
$ws.Range("P28").Value = 0.1932063174578422
$ws.Range("Q28").Value = 28.28353995926667
$ws.Range("R28").Value = 254.5518596334
$ws.Range("S28").Value = 0.009137274231579029
$ws.Range("T28").Value = 0.009137274231579028

$ws.Range("E29").Value = 2
$ws.Range("F29").Value = 0.6666666666666666
$ws.Range("G29").Value = 0.41099
$ws.Range("H29").Value = 1.23297
$ws.Range("I29").Value = 0.0472928336495663
$ws.Range("J29").Value = 0.0472928336495663
$ws.Range("M29").Value = 11.78107633333333
$ws.Range("N29").Value = 35.343229
$ws.Range("O29").Value = 0.03307529931894448
$ws.Range("P29").Value = 0.03307529931894448
$ws.Range("Q29").Value = 4.841904562236667
$ws.Range("R29").Value = 43.57714106013
$ws.Range("S29").Value = 0.001564224628600455
$ws.Range("T29").Value = 0.001564224628600455

$ws.Range("E30").Value = 2
$ws.Range("F30").Value = 0.6666666666666666
$ws.Range("G30").Value = 0.41099
$ws.Range("H30").Value = 1.23297
$ws.Range("I30").Value = 0.0472928336495663
$ws.Range("J30").Value = 0.0472928336495663
$ws.Range("M30").Value = 139.820737
$ws.Range("N30").Value = 419.462211
$ws.Range("O30").Value = 0.3925458588351179
$ws.Range("P30").Value = 0.3925458588351179
$ws.Range("Q30").Value = 57.46492469963
$ws.Range("R30").Value = 517.18432229667
$ws.Range("S30").Value = 0.01856460600171537
$ws.Range("T30").Value = 0.01856460600171537

$ws.Range("E31").Value = 2
$ws.Range("F31").Value = 0.6666666666666666
$ws.Range("G31").Value = 0.41099
$ws.Range("H31").Value = 1.23297
$ws.Range("I31").Value = 0.0472928336495663
$ws.Range("J31").Value = 0.0472928336495663
$ws.Range("M31").Value = 40.90925733333334
$ws.Range("N31").Value = 122.727772
$ws.Range("O31").Value = 0.1148524882558742
$ws.Range("P31").Value = 0.1148524882558742
$ws.Range("Q31").Value = 16.81329567142667
$ws.Range("R31").Value = 151.31966104284
$ws.Range("S31").Value = 0.005431699621323827
$ws.Range("T31").Value = 0.005431699621323826

$ws.Range("G32").Value = 3.718746666666667
$ws.Range("H32").Value = 11.15624
$ws.Range("I32").Value = 0.4279181184251341
$ws.Range("J32").Value = 0.4279181184251342
$ws.Range("M32").Value = 42.04602466666667
$ws.Range("N32").Value = 126.138074
$ws.Range("O32").Value = 0.1180439555498783
$ws.Range("P32").Value = 0.1180439555498783
$ws.Range("Q32").Value = 156.3585140757511
$ws.Range("R32").Value = 1407.22662668176
$ws.Range("S32").Value = 0.05051314735036407
$ws.Range("T32").Value = 0.05051314735036408

$ws.Range("G33").Value = 3.718746666666667
$ws.Range("H33").Value = 11.15624
$ws.Range("I33").Value = 0.4279181184251341
$ws.Range("J33").Value = 0.4279181184251342
$ws.Range("O33").Value = 0.1482760805823429
$ws.Range("P33").Value = 0.1482760805823429
$ws.Range("Q33").Value = 196.4033442020267
$ws.Range("R33").Value = 1767.63009781824
$ws.Range("S33").Value = 0.06345002141024975
$ws.Range("T33").Value = 0.06345002141024976

$ws.Range("G34").Value = 3.718746666666667
$ws.Range("H34").Value = 11.15624
$ws.Range("I34").Value = 0.4279181184251341
$ws.Range("J34").Value = 0.4279181184251342
$ws.Range("M34").Value = 68.81807333333334
$ws.Range("N34").Value = 206.45422
$ws.Range("O34").Value = 0.1932063174578422
$ws.Range("P34").Value = 0.1932063174578422
$ws.Range("Q34").Value = 255.9169808147556
$ws.Range("R34").Value = 2303.2528273328
$ws.Range("S34").Value = 0.08267648383440897
$ws.Range("T34").Value = 0.08267648383440897

$ws.Range("G35").Value = 3.718746666666667
$ws.Range("H35").Value = 11.15624
$ws.Range("I35").Value = 0.4279181184251341
$ws.Range("J35").Value = 0.4279181184251342
$ws.Range("M35").Value = 11.78107633333333
$ws.Range("N35").Value = 35.343229
$ws.Range("O35").Value = 0.03307529931894448
$ws.Range("P35").Value = 0.03307529931894448
$ws.Range("Q35").Value = 43.81083834432889
$ws.Range("R35").Value = 394.29754509896
$ws.Range("S35").Value = 0.01415351985091084
$ws.Range("T35").Value = 0.01415351985091084

$ws.Range("G36").Value = 3.718746666666667
$ws.Range("H36").Value = 11.15624
$ws.Range("I36").Value = 0.4279181184251341
$ws.Range("J36").Value = 0.4279181184251342
$ws.Range("M36").Value = 139.820737
$ws.Range("N36").Value = 419.462211
$ws.Range("O36").Value = 0.3925458588351179
$ws.Range("P36").Value = 0.3925458588351179
$ws.Range("Q36").Value = 519.9578996496267
$ws.Range("R36").Value = 4679.62109684664
$ws.Range("S36").Value = 0.1679774853083019
$ws.Range("T36").Value = 0.167977485308302

$ws.Range("G37").Value = 3.718746666666667
$ws.Range("H37").Value = 11.15624
$ws.Range("I37").Value = 0.4279181184251341
$ws.Range("J37").Value = 0.4279181184251342
$ws.Range("M37").Value = 40.90925733333334
$ws.Range("N37").Value = 122.727772
$ws.Range("O37").Value = 0.1148524882558742
$ws.Range("P37").Value = 0.1148524882558742
$ws.Range("Q37").Value = 152.1311643441422
$ws.Range("R37").Value = 1369.18047909728
$ws.Range("S37").Value = 0.04914746067089851
$ws.Range("T37").Value = 0.04914746067089851
